$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INCO")

for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 10).Value = "Cambio"
    $ws.Cells.Item($row, 11).Value = "Sin equipos"
}
